{"js": "// Applies the Xhosa-translation edits described in the diff.\n// Each change replaces an exact run of text with its translated\n// counterpart, leaving paragraph/run formatting untouched.\n\nconst body = context.document.body;\n\nasync function replaceOnce(oldText, newText) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n  results.items[0].insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 1) Invitation sentence translated into Xhosa.\nawait replaceOnce(\n  \"You've been invited to an interview with a member of our research team because you\\u2019re part of our study. \",\n  \"Umenywe kudliwano-ndlebe nelungu leqela lethu lophando kuba uyinxalenye yophononongo lwethu. \"\n);\n\n// 2) Audio-recording sentence translated into Xhosa.\nawait replaceOnce(\n  \"We will audio record the interview to help us remember the discussion and later write down what was said. \",\n  \"Sizakurekhoda oludliwano-ndlebe ukusinceda sikhumbule ebesixoxe ngako kwaye kamva sikubhale phantsi obekuthethiwe. \"\n);\n\n// 3) \"However, we will be\" / \"interviewing many\" / \" nurses and clinic staff...\" split\n//    across three runs \u2014 replace each run's text independently.\nawait replaceOnce(\n  \"However, we will be \",\n  \"Nangona kunjalo, siza kuba \"\n);\n\nawait replaceOnce(\n  \"interviewing many\",\n  \"udliwano-ndlebe nabongi\"\n);\n\nawait replaceOnce(\n  \" nurses and clinic staff from at least 7 clinics across two sites in the Western Cape, and we will be carefully watching out for any details you share which may accidentally reveal your identity. \",\n  \" abaninzi kunye nabasebenzi basekliniki ubuncinane kwiikliniki ezisixhenxe kwiindawo ezibini eNtshona Koloni, kwaye siyakube sijonge ngononophelo kuzo naziphi na iinkcukacha owabelana ngazo ezinokutyhila ubuwena ngempazamo. \"\n);\n\n// 4) \"Interview recordings will be deleted after transcription.\" translated.\nawait replaceOnce(\n  \"Interview recordings will be deleted after transcription. \",\n  \"Urekhodo lodliwano-ndlebe luya kucinywa emva kokukhutshelwa. \"\n);\n\n// 5) Consent checkbox item translated into Xhosa.\nawait replaceOnce(\n  \"I am okay with the interview being recorded. I know the recordings will be used for research.\",\n  \"Ndikulungele ukuba nodliwano-ndlebe olurekhodiweyo. Ndiyayazi iirekhodingi zizakusetyenziswa kuphando.\"\n);\n\n// 6) Signature label translated into Xhosa.\nawait replaceOnce(\n  \"Signature of Interviewee\",\n  \"Umtyikityo womntu owenziwa udliwano-ndlebe\"\n);\n", "ps1": "# Applies the Xhosa-translation edits described in the diff.\n# Each replacement targets one exact run of text and swaps it for its\n# translated counterpart, leaving all other formatting untouched.\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($find, $replace) {\n    $range = $d.Content\n    $range.Find.ClearFormatting()\n    $range.Find.Replacement.ClearFormatting()\n    $range.Find.Text = $find\n    $range.Find.Replacement.Text = $replace\n    # NOTE: MatchCase=$true combined with a trailing space in $find is\n    # mishandled by this host's Find engine (it reports no match even\n    # though the text is present), so we search case-insensitively here.\n    # Every pattern below is unique in the document either way, so this\n    # does not risk an unintended match.\n    $range.Find.Execute($find, $false, $false, $false, $false, $false, $true, 1, $false, $replace, 2)\n}\n\n# 1) Invitation sentence translated into Xhosa.\nReplace-Text \"You've been invited to an interview with a member of our research team because you\u2019re part of our study. \" \"Umenywe kudliwano-ndlebe nelungu leqela lethu lophando kuba uyinxalenye yophononongo lwethu. \"\n\n# 2) Audio-recording sentence translated into Xhosa.\nReplace-Text \"We will audio record the interview to help us remember the discussion and later write down what was said. \" \"Sizakurekhoda oludliwano-ndlebe ukusinceda sikhumbule ebesixoxe ngako kwaye kamva sikubhale phantsi obekuthethiwe. \"\n\n# 3) \"However, we will be\" / \"interviewing many\" / \" nurses and clinic staff...\" are\n#    three separate runs in the source paragraph \u2014 replace each run's text in turn.\nReplace-Text \"However, we will be \" \"Nangona kunjalo, siza kuba \"\nReplace-Text \"interviewing many\" \"udliwano-ndlebe nabongi\"\nReplace-Text \" nurses and clinic staff from at least 7 clinics across two sites in the Western Cape, and we will be carefully watching out for any details you share which may accidentally reveal your identity. \" \" abaninzi kunye nabasebenzi basekliniki ubuncinane kwiikliniki ezisixhenxe kwiindawo ezibini eNtshona Koloni, kwaye siyakube sijonge ngononophelo kuzo naziphi na iinkcukacha owabelana ngazo ezinokutyhila ubuwena ngempazamo. \"\n\n# 4) \"Interview recordings will be deleted after transcription.\" translated.\nReplace-Text \"Interview recordings will be deleted after transcription. \" \"Urekhodo lodliwano-ndlebe luya kucinywa emva kokukhutshelwa. \"\n\n# 5) Consent checkbox item translated into Xhosa.\nReplace-Text \"I am okay with the interview being recorded. I know the recordings will be used for research.\" \"Ndikulungele ukuba nodliwano-ndlebe olurekhodiweyo. Ndiyayazi iirekhodingi zizakusetyenziswa kuphando.\"\n\n# 6) Signature label translated into Xhosa.\nReplace-Text \"Signature of Interviewee\" \"Umtyikityo womntu owenziwa udliwano-ndlebe\"\n"}
